# Applies crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook

# Temporary staging sheet used to force numeric-looking price strings
# to remain plain text (avoids Excel auto-converting "0.619" etc. to a number).
$tmp = $wb.Worksheets.Add()
$tmp.Range("A1").Value = "'244.67"
$tmp.Range("A2").Value = "'0.619"
$tmp.Range("A3").Value = "'59.39"
$tmp.Range("A4").Value = "'0.375"
$tmp.Range("A5").Value = "'0.0814"
$tmp.Range("A6").Value = "'22.30"
$tmp.Range("A7").Value = "'13.72"
$tmp.Range("A8").Value = "'70.05"
$tmp.Range("A9").Value = "'229.32"
$tmp.Range("A10").Value = "'9.25"
$tmp.Range("A11").Value = "'0.141"
$tmp.Range("A12").Value = "'160.11"
$tmp.Range("A13").Value = "'19.40"
$tmp.Range("A14").Value = "'0.0619"
$tmp.Range("A15").Value = "'1.00"
$tmp.Range("A16").Value = "'2.27"
$tmp.Range("A17").Value = "'3.41"
$tmp.Range("A18").Value = "'0.0988"
$tmp.Range("A19").Value = "'16.16"
$tmp.Range("A20").Value = "'87.92"
$tmp.Range("A21").Value = "'2.84"
$tmp.Range("A22").Value = "'43.86"

# Re-fetch Sheet1 now that a new sheet has been added (keeps reference fresh).
$ws = $wb.Worksheets.Item("Sheet1")

# --- Plain text / safe updates (Coin, Link, Volume%, and non-numeric-looking Price cells) ---
$ws.Range("D2").Value = "36.618.26"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.964.54"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  +2.10%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("E12").Value = "  +3.60%  "
$ws.Range("D13").Value = "2.252.66"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "1.966.61"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "36.558.39"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "0.0₃0859"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("E25").Value = "  +3.13%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E27").Value = "  +11.57%  "
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("B35").Value = "BinanceUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E36").Value = "  +6.28%  "
$ws.Range("E37").Value = "  +12.76%  "
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("D45").Value = "1.359.69"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").Value = "2.143.75"
$ws.Range("E50").Value = "  +0.95%  "
$ws.Range("E51").Value = "  -3.52%  "

# --- Price cells that look like plain numbers: copy text value in from staging sheet ---
$tmp.Range("A1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$tmp.Range("A2").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$tmp.Range("A3").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$tmp.Range("A4").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$tmp.Range("A5").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$tmp.Range("A6").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$tmp.Range("A7").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$tmp.Range("A8").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$tmp.Range("A9").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$tmp.Range("A10").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$tmp.Range("A11").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$tmp.Range("A12").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$tmp.Range("A13").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$tmp.Range("A14").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$tmp.Range("A15").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$tmp.Range("A16").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$tmp.Range("A17").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$tmp.Range("A18").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$tmp.Range("A19").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$tmp.Range("A20").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$tmp.Range("A21").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$tmp.Range("A22").Copy()
$ws.Range("D51").PasteSpecial(-4163)

# Remove the staging sheet; clear clipboard/marching-ants selection state.
[void]$tmp.Delete()
$excel.CutCopyMode = 0

Write-Host "Applied cryptos list update"
